# "updated phase2 practice removed extra cols"
# Column E was an unused/extra column; the "Time" column (currently F)
# is moved left into E so the data is contiguous (A:E) instead of A:D + F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move F1:F5 ("Time" header + values) one column left, into E1:E5.
$ws.Range("F1:F5").Cut($ws.Range("E1:E5"))

# Restore the (now out-of-data) selection that was left on the sheet.
$ws.Range("F11").Select()
